$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of the last existing data row (357) down into
# the new rows (358:366) so the new "A" column cells keep the same date
# number-format / style index as the rest of the column.
$srcRow = $ws.Range("A357:D357")
$dstRow = $ws.Range("A358:D366")
$srcRow.Copy($dstRow)

# New daily data (date serial, nuovi pos., somma mobile 7gg., somma mobile
# 7gg. per 100mila abitanti) bringing the report up to 1/09/2021.
$data = @(
    @(44432, 0, 1, 14.22475106685633),
    @(44433, 0, 1, 14.22475106685633),
    @(44434, 0, 1, 14.22475106685633),
    @(44435, 2, 3, 42.67425320056899),
    @(44436, 0, 3, 42.67425320056899),
    @(44437, 1, 3, 42.67425320056899),
    @(44438, 0, 3, 42.67425320056899),
    @(44439, 2, 5, 71.12375533428165),
    @(44440, 1, 6, 85.34850640113798)
)

$r = 358
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
